$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the model-fit data for row 30 (age 75, cause index for A30/B30/C30 already present)
$ws.Range("D30").Value = 134331.96
$ws.Range("E30").Value = 134337.22
$ws.Range("F30").Formula = "=IF(D30<E30,""non-pw"",""pw"")"

# Update the view: scroll so row 19 is the top-left row, and select F30:F31 with F31 active
$ws.Range("F31").Select()
$ws.Range("F30:F31").Select()
$excel.ActiveWindow.ScrollRow = 19
